$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Remove the two obsolete "2507" period rows (DIANA SOFIA RODRIGUEZ
# BARRIOS and MARELIS MADRID TOUS) - rows 18 and 19. This shifts the
# remaining rows (old 20, 21 and the signature block) up by two.
# ------------------------------------------------------------------
$ws.Range("B18:J19").EntireRow.Delete()

# After the delete:
#   row 18 -> KAOLIS BOSSA BALLESTERO / 2508 (unchanged, keep as-is)
#   row 19 -> MARELIS MADRID TOUS / 2508 (last row, bottom-bordered style)
# Insert a brand-new row 20 below it (for the new KAOLIS / 2509 entry),
# pushing the signature block down, and give it the same border
# formatting as row 19 (the bottom-bordered "last row" look).
$ws.Rows.Item(20).Insert()
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# Row 19 is no longer the last table row, so restyle it like the
# regular (non-bordered) data rows 16-18, then overwrite its content
# with the new worker: NURIS DEL CARMEN SARMIENTO GUETE.
$ws.Range("B16:J16").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45476981"
$ws.Range("D19").Value = "NURIS DEL CARMEN SARMIENTO GUETE"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 77334
$ws.Range("G19").Value = 2000000

# New last row (20): KAOLIS BOSSA BALLESTERO again, now for period 2509.
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143375582"
$ws.Range("D20").Value = "KAOLIS BOSSA BALLESTERO"
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# ------------------------------------------------------------------
# Header summary updates
# ------------------------------------------------------------------
$ws.Range("E11").Value = 223155
$ws.Range("C13").Value = 4

# ------------------------------------------------------------------
# Column D got a bit wider to fit the new, longer worker name.
# ------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 36.6328125
